$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

# Update Version and Date values
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" row before row 11 ("Description"), pushing existing rows down
$meta.Rows.Item(11).Insert()
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# Match formatting of the other data rows (style was lost on insert)
$meta.Range("A12:B12").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""
